$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Ntrk3"
$ws.Cells.Item(2, 3).Value = "Ptprf"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.1511493333333333
$ws.Cells.Item(2, 8).Value = 0.453448
$ws.Cells.Item(2, 9).Value = 0.7495144539818079
$ws.Cells.Item(2, 10).Value = 0.7495144539818078
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.101297
$ws.Cells.Item(2, 14).Value = 0.303891
$ws.Cells.Item(2, 15).Value = 0.009886476207827108
$ws.Cells.Item(2, 16).Value = 0.009886476207827108
$ws.Cells.Item(2, 17).Value = 0.01531097401866667
$ws.Cells.Item(2, 18).Value = 0.137798766168
$ws.Cells.Item(2, 19).Value = 0.00741005681671367
$ws.Cells.Item(2, 20).Value = 0.007410056816713669

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Ntrk3"
$ws.Cells.Item(3, 3).Value = "Ptprf"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.1511493333333333
$ws.Cells.Item(3, 8).Value = 0.453448
$ws.Cells.Item(3, 9).Value = 0.7495144539818079
$ws.Cells.Item(3, 10).Value = 0.7495144539818078
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 5.829902333333333
$ws.Cells.Item(3, 14).Value = 17.489707
$ws.Cells.Item(3, 15).Value = 0.5689920798489169
$ws.Cells.Item(3, 16).Value = 0.5689920798489169
$ws.Cells.Item(3, 17).Value = 0.8811858510817777
$ws.Cells.Item(3, 18).Value = 7.930672659735999
$ws.Cells.Item(3, 19).Value = 0.4264677880479342
$ws.Cells.Item(3, 20).Value = 0.4264677880479341

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Ntrk3"
$ws.Cells.Item(4, 3).Value = "Ptprf"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.1511493333333333
$ws.Cells.Item(4, 8).Value = 0.453448
$ws.Cells.Item(4, 9).Value = 0.7495144539818079
$ws.Cells.Item(4, 10).Value = 0.7495144539818078
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.06457233333333333
$ws.Cells.Item(4, 14).Value = 0.193717
$ws.Cells.Item(4, 15).Value = 0.006302188980758376
$ws.Cells.Item(4, 16).Value = 0.006302188980758376
$ws.Cells.Item(4, 17).Value = 0.00976006513511111
$ws.Cells.Item(4, 18).Value = 0.087840586216
$ws.Cells.Item(4, 19).Value = 0.004723581732803281
$ws.Cells.Item(4, 20).Value = 0.00472358173280328

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ntrk3"
$ws.Cells.Item(5, 3).Value = "Ptprf"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.1511493333333333
$ws.Cells.Item(5, 8).Value = 0.453448
$ws.Cells.Item(5, 9).Value = 0.7495144539818079
$ws.Cells.Item(5, 10).Value = 0.7495144539818078
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.250245
$ws.Cells.Item(5, 14).Value = 12.750735
$ws.Cells.Item(5, 15).Value = 0.4148192549624977
$ws.Cells.Item(5, 16).Value = 0.4148192549624976
$ws.Cells.Item(5, 17).Value = 0.6424216982533334
$ws.Cells.Item(5, 18).Value = 5.78179528428
$ws.Cells.Item(5, 19).Value = 0.3109130273843568
$ws.Cells.Item(5, 20).Value = 0.3109130273843567

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Ntrk3"
$ws.Cells.Item(6, 3).Value = "Ptprf"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.05051366666666667
$ws.Cells.Item(6, 8).Value = 0.151541
$ws.Cells.Item(6, 9).Value = 0.2504855460181921
$ws.Cells.Item(6, 10).Value = 0.2504855460181921
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.101297
$ws.Cells.Item(6, 14).Value = 0.303891
$ws.Cells.Item(6, 15).Value = 0.009886476207827108
$ws.Cells.Item(6, 16).Value = 0.009886476207827108
$ws.Cells.Item(6, 17).Value = 0.005116882892333335
$ws.Cells.Item(6, 18).Value = 0.046051946031
$ws.Cells.Item(6, 19).Value = 0.002476419391113438
$ws.Cells.Item(6, 20).Value = 0.002476419391113438

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Ntrk3"
$ws.Cells.Item(7, 3).Value = "Ptprf"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.05051366666666667
$ws.Cells.Item(7, 8).Value = 0.151541
$ws.Cells.Item(7, 9).Value = 0.2504855460181921
$ws.Cells.Item(7, 10).Value = 0.2504855460181921
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.829902333333333
$ws.Cells.Item(7, 14).Value = 17.489707
$ws.Cells.Item(7, 15).Value = 0.5689920798489169
$ws.Cells.Item(7, 16).Value = 0.5689920798489169
$ws.Cells.Item(7, 17).Value = 0.2944897431652223
$ws.Cells.Item(7, 18).Value = 2.650407688487
$ws.Cells.Item(7, 19).Value = 0.1425242918009827
$ws.Cells.Item(7, 20).Value = 0.1425242918009827

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ntrk3"
$ws.Cells.Item(8, 3).Value = "Ptprf"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.05051366666666667
$ws.Cells.Item(8, 8).Value = 0.151541
$ws.Cells.Item(8, 9).Value = 0.2504855460181921
$ws.Cells.Item(8, 10).Value = 0.2504855460181921
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.06457233333333333
$ws.Cells.Item(8, 14).Value = 0.193717
$ws.Cells.Item(8, 15).Value = 0.006302188980758376
$ws.Cells.Item(8, 16).Value = 0.006302188980758376
$ws.Cells.Item(8, 17).Value = 0.003261785321888889
$ws.Cells.Item(8, 18).Value = 0.029356067897
$ws.Cells.Item(8, 19).Value = 0.001578607247955095
$ws.Cells.Item(8, 20).Value = 0.001578607247955095

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ntrk3"
$ws.Cells.Item(9, 3).Value = "Ptprf"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.05051366666666667
$ws.Cells.Item(9, 8).Value = 0.151541
$ws.Cells.Item(9, 9).Value = 0.2504855460181921
$ws.Cells.Item(9, 10).Value = 0.2504855460181921
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 4.250245
$ws.Cells.Item(9, 14).Value = 12.750735
$ws.Cells.Item(9, 15).Value = 0.4148192549624977
$ws.Cells.Item(9, 16).Value = 0.4148192549624976
$ws.Cells.Item(9, 17).Value = 0.2146954591816667
$ws.Cells.Item(9, 18).Value = 1.932259132635
$ws.Cells.Item(9, 19).Value = 0.1039062275781409
$ws.Cells.Item(9, 20).Value = 0.1039062275781408
